$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 31.824752
$ws.Cells.Item(2, 8).Value = 95.47425600000001
$ws.Cells.Item(2, 9).Value = 0.886907633630525
$ws.Cells.Item(2, 10).Value = 0.886907633630525
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.042868
$ws.Cells.Item(2, 14).Value = 0.128604
$ws.Cells.Item(2, 15).Value = 0.03014606792405771
$ws.Cells.Item(2, 16).Value = 0.03014606792405771
$ws.Cells.Item(2, 17).Value = 1.364263468736
$ws.Cells.Item(2, 18).Value = 12.278371218624
$ws.Cells.Item(2, 19).Value = 0.0267367777657911
$ws.Cells.Item(2, 20).Value = 0.0267367777657911

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 31.824752
$ws.Cells.Item(3, 8).Value = 95.47425600000001
$ws.Cells.Item(3, 9).Value = 0.886907633630525
$ws.Cells.Item(3, 10).Value = 0.886907633630525
$ws.Cells.Item(3, 15).Value = 0.2718481285523376
$ws.Cells.Item(3, 16).Value = 0.2718481285523376
$ws.Cells.Item(3, 17).Value = 12.30251559714134
$ws.Cells.Item(3, 18).Value = 110.722640374272
$ws.Cells.Item(3, 19).Value = 0.2411041804012405
$ws.Cells.Item(3, 20).Value = 0.2411041804012405

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 31.824752
$ws.Cells.Item(4, 8).Value = 95.47425600000001
$ws.Cells.Item(4, 9).Value = 0.886907633630525
$ws.Cells.Item(4, 10).Value = 0.886907633630525
$ws.Cells.Item(4, 13).Value = 0.9839956666666666
$ws.Cells.Item(4, 14).Value = 2.951987
$ws.Cells.Item(4, 15).Value = 0.69197537100662
$ws.Cells.Item(4, 16).Value = 0.69197537100662
$ws.Cells.Item(4, 17).Value = 31.31541806074134
$ws.Cells.Item(4, 18).Value = 281.838762546672
$ws.Cells.Item(4, 19).Value = 0.6137182388300859
$ws.Cells.Item(4, 20).Value = 0.6137182388300859

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 31.824752
$ws.Cells.Item(5, 8).Value = 95.47425600000001
$ws.Cells.Item(5, 9).Value = 0.886907633630525
$ws.Cells.Item(5, 10).Value = 0.886907633630525
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.008575333333333332
$ws.Cells.Item(5, 14).Value = 0.025726
$ws.Cells.Item(5, 15).Value = 0.006030432516984765
$ws.Cells.Item(5, 16).Value = 0.006030432516984765
$ws.Cells.Item(5, 17).Value = 0.2729078566506667
$ws.Cells.Item(5, 18).Value = 2.456170709856
$ws.Cells.Item(5, 19).Value = 0.005348436633407528
$ws.Cells.Item(5, 20).Value = 0.005348436633407529

$ws.Cells.Item(6, 9).Value = 0.06502043684278042
$ws.Cells.Item(6, 10).Value = 0.06502043684278042
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.042868
$ws.Cells.Item(6, 14).Value = 0.128604
$ws.Cells.Item(6, 15).Value = 0.03014606792405771
$ws.Cells.Item(6, 16).Value = 0.03014606792405771
$ws.Cells.Item(6, 17).Value = 0.100016059556
$ws.Cells.Item(6, 18).Value = 0.900144536004
$ws.Cells.Item(6, 19).Value = 0.001960110505514363
$ws.Cells.Item(6, 20).Value = 0.001960110505514363

$ws.Cells.Item(7, 9).Value = 0.06502043684278042
$ws.Cells.Item(7, 10).Value = 0.06502043684278042
$ws.Cells.Item(7, 15).Value = 0.2718481285523376
$ws.Cells.Item(7, 16).Value = 0.2718481285523376
$ws.Cells.Item(7, 19).Value = 0.01767568407336532
$ws.Cells.Item(7, 20).Value = 0.01767568407336532

$ws.Cells.Item(8, 9).Value = 0.06502043684278042
$ws.Cells.Item(8, 10).Value = 0.06502043684278042
$ws.Cells.Item(8, 13).Value = 0.9839956666666666
$ws.Cells.Item(8, 14).Value = 2.951987
$ws.Cells.Item(8, 15).Value = 0.69197537100662
$ws.Cells.Item(8, 16).Value = 0.69197537100662
$ws.Cells.Item(8, 17).Value = 2.295777017826333
$ws.Cells.Item(8, 18).Value = 20.661993160437
$ws.Cells.Item(8, 19).Value = 0.04499254090729548
$ws.Cells.Item(8, 20).Value = 0.04499254090729548

$ws.Cells.Item(9, 9).Value = 0.06502043684278042
$ws.Cells.Item(9, 10).Value = 0.06502043684278042
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.008575333333333332
$ws.Cells.Item(9, 14).Value = 0.025726
$ws.Cells.Item(9, 15).Value = 0.006030432516984765
$ws.Cells.Item(9, 16).Value = 0.006030432516984765
$ws.Cells.Item(9, 17).Value = 0.02000725598066666
$ws.Cells.Item(9, 18).Value = 0.180065303826
$ws.Cells.Item(9, 19).Value = 0.0003921013566052572
$ws.Cells.Item(9, 20).Value = 0.0003921013566052573

$ws.Cells.Item(10, 7).Value = 1.696588
$ws.Cells.Item(10, 8).Value = 5.089764000000001
$ws.Cells.Item(10, 9).Value = 0.04728133775640876
$ws.Cells.Item(10, 10).Value = 0.04728133775640876
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.042868
$ws.Cells.Item(10, 14).Value = 0.128604
$ws.Cells.Item(10, 15).Value = 0.03014606792405771
$ws.Cells.Item(10, 16).Value = 0.03014606792405771
$ws.Cells.Item(10, 17).Value = 0.07272933438400001
$ws.Cells.Item(10, 18).Value = 0.654564009456
$ws.Cells.Item(10, 19).Value = 0.001425346419545013
$ws.Cells.Item(10, 20).Value = 0.001425346419545013

$ws.Cells.Item(11, 7).Value = 1.696588
$ws.Cells.Item(11, 8).Value = 5.089764000000001
$ws.Cells.Item(11, 9).Value = 0.04728133775640876
$ws.Cells.Item(11, 10).Value = 0.04728133775640876
$ws.Cells.Item(11, 15).Value = 0.2718481285523376
$ws.Cells.Item(11, 16).Value = 0.2718481285523376
$ws.Cells.Item(11, 17).Value = 0.6558511542186668
$ws.Cells.Item(11, 18).Value = 5.902660387968001
$ws.Cells.Item(11, 19).Value = 0.0128533431845307
$ws.Cells.Item(11, 20).Value = 0.0128533431845307

$ws.Cells.Item(12, 7).Value = 1.696588
$ws.Cells.Item(12, 8).Value = 5.089764000000001
$ws.Cells.Item(12, 9).Value = 0.04728133775640876
$ws.Cells.Item(12, 10).Value = 0.04728133775640876
$ws.Cells.Item(12, 13).Value = 0.9839956666666666
$ws.Cells.Item(12, 14).Value = 2.951987
$ws.Cells.Item(12, 15).Value = 0.69197537100662
$ws.Cells.Item(12, 16).Value = 0.69197537100662
$ws.Cells.Item(12, 17).Value = 1.669435240118667
$ws.Cells.Item(12, 18).Value = 15.024917161068
$ws.Cells.Item(12, 19).Value = 0.03271752123568026
$ws.Cells.Item(12, 20).Value = 0.03271752123568026

$ws.Cells.Item(13, 7).Value = 1.696588
$ws.Cells.Item(13, 8).Value = 5.089764000000001
$ws.Cells.Item(13, 9).Value = 0.04728133775640876
$ws.Cells.Item(13, 10).Value = 0.04728133775640876
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.008575333333333332
$ws.Cells.Item(13, 14).Value = 0.025726
$ws.Cells.Item(13, 15).Value = 0.006030432516984765
$ws.Cells.Item(13, 16).Value = 0.006030432516984765
$ws.Cells.Item(13, 17).Value = 0.01454880762933333
$ws.Cells.Item(13, 18).Value = 0.130939268664
$ws.Cells.Item(13, 19).Value = 0.0002851269166527868
$ws.Cells.Item(13, 20).Value = 0.0002851269166527869

$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.02836866666666667
$ws.Cells.Item(14, 8).Value = 0.085106
$ws.Cells.Item(14, 9).Value = 0.0007905917702857979
$ws.Cells.Item(14, 10).Value = 0.0007905917702857978
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.042868
$ws.Cells.Item(14, 14).Value = 0.128604
$ws.Cells.Item(14, 15).Value = 0.03014606792405771
$ws.Cells.Item(14, 16).Value = 0.03014606792405771
$ws.Cells.Item(14, 17).Value = 0.001216108002666667
$ws.Cells.Item(14, 18).Value = 0.010944972024
$ws.Cells.Item(14, 19).Value = 0.00002383323320723669
$ws.Cells.Item(14, 20).Value = 0.00002383323320723669

$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.02836866666666667
$ws.Cells.Item(15, 8).Value = 0.085106
$ws.Cells.Item(15, 9).Value = 0.0007905917702857979
$ws.Cells.Item(15, 10).Value = 0.0007905917702857978
$ws.Cells.Item(15, 15).Value = 0.2718481285523376
$ws.Cells.Item(15, 16).Value = 0.2718481285523376
$ws.Cells.Item(15, 17).Value = 0.01096649438577778
$ws.Cells.Item(15, 18).Value = 0.098698449472
$ws.Cells.Item(15, 19).Value = 0.0002149208932010737
$ws.Cells.Item(15, 20).Value = 0.0002149208932010737

$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.02836866666666667
$ws.Cells.Item(16, 8).Value = 0.085106
$ws.Cells.Item(16, 9).Value = 0.0007905917702857979
$ws.Cells.Item(16, 10).Value = 0.0007905917702857978
$ws.Cells.Item(16, 13).Value = 0.9839956666666666
$ws.Cells.Item(16, 14).Value = 2.951987
$ws.Cells.Item(16, 15).Value = 0.69197537100662
$ws.Cells.Item(16, 16).Value = 0.69197537100662
$ws.Cells.Item(16, 17).Value = 0.02791464506911111
$ws.Cells.Item(16, 18).Value = 0.251231805622
$ws.Cells.Item(16, 19).Value = 0.0005470700335582956
$ws.Cells.Item(16, 20).Value = 0.0005470700335582954

$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.02836866666666667
$ws.Cells.Item(17, 8).Value = 0.085106
$ws.Cells.Item(17, 9).Value = 0.0007905917702857979
$ws.Cells.Item(17, 10).Value = 0.0007905917702857978
$ws.Cells.Item(17, 11).Value = 1
$ws.Cells.Item(17, 12).Value = 0.3333333333333333
$ws.Cells.Item(17, 13).Value = 0.008575333333333332
$ws.Cells.Item(17, 14).Value = 0.025726
$ws.Cells.Item(17, 15).Value = 0.006030432516984765
$ws.Cells.Item(17, 16).Value = 0.006030432516984765
$ws.Cells.Item(17, 17).Value = 0.0002432707728888889
$ws.Cells.Item(17, 18).Value = 0.002189436956
$ws.Cells.Item(17, 19).Value = 0.000004767610319192025
$ws.Cells.Item(17, 20).Value = 0.000004767610319192025
